$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.466.05"
$ws.Range("D2").Style = $ws.Range("D6").Style
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "'1.819.01"
$ws.Range("D3").Style = $ws.Range("D6").Style
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = $ws.Range("D6").Style
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'316.55"
$ws.Range("D5").Style = $ws.Range("D6").Style
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "'0.5154"
$ws.Range("D7").Style = $ws.Range("D6").Style
$ws.Range("E7").Value = "  -3.37%  "

$ws.Range("D8").Value = "'0.3887"
$ws.Range("D8").Style = $ws.Range("D6").Style
$ws.Range("E8").Value = "  -2.58%  "

$ws.Range("D9").Value = "'0.08482"
$ws.Range("D9").Style = $ws.Range("D6").Style
$ws.Range("E9").Value = "  +8.99%  "

$ws.Range("D10").Value = "'41.82"
$ws.Range("D10").Style = $ws.Range("D6").Style
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").Value = "'6.440"
$ws.Range("D12").Style = $ws.Range("D6").Style
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").Value = "'1.003"
$ws.Range("D14").Style = $ws.Range("D6").Style
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").Value = "'7.501"
$ws.Range("D15").Style = $ws.Range("D6").Style
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").Value = "'1.822.42"
$ws.Range("D16").Style = $ws.Range("D6").Style
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").Value = "'0.00001141"
$ws.Range("D17").Style = $ws.Range("D6").Style
$ws.Range("E17").Value = "  +4.39%  "

$ws.Range("D18").Value = "'92.77"
$ws.Range("D18").Style = $ws.Range("D6").Style
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").Value = "'0.06654"
$ws.Range("D19").Style = $ws.Range("D6").Style
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "'17.73"
$ws.Range("D20").Style = $ws.Range("D6").Style
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = $ws.Range("D6").Style
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "'6.080"
$ws.Range("D22").Style = $ws.Range("D6").Style
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").Value = "'28.504.86"
$ws.Range("D23").Style = $ws.Range("D6").Style
$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("D24").Value = "'11.41"
$ws.Range("D24").Style = $ws.Range("D6").Style
$ws.Range("E24").Value = "  +1.71%  "

$ws.Range("D25").Value = "'2.274"
$ws.Range("D25").Style = $ws.Range("D6").Style
$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("D26").Value = "'21.02"
$ws.Range("D26").Style = $ws.Range("D6").Style
$ws.Range("E26").Value = "  +0.82%  "

$ws.Range("D27").Value = "'159.29"
$ws.Range("D27").Style = $ws.Range("D6").Style
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("D28").Value = "'2.030.34"
$ws.Range("D28").Style = $ws.Range("D6").Style
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "'2.405"
$ws.Range("D29").Style = $ws.Range("D6").Style
$ws.Range("E29").Value = "  -1.43%  "

$ws.Range("D30").Value = "'125.68"
$ws.Range("D30").Style = $ws.Range("D6").Style
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("E31").Value = "  -4.07%  "

$ws.Range("D32").Value = "'1.092"
$ws.Range("D32").Style = $ws.Range("D6").Style
$ws.Range("E32").Value = "  -5.84%  "

$ws.Range("D33").Value = "'5.721"
$ws.Range("D33").Style = $ws.Range("D6").Style
$ws.Range("E33").Value = "  -0.83%  "

$ws.Range("D34").Value = "'0.07458"
$ws.Range("D34").Style = $ws.Range("D6").Style
$ws.Range("E34").Value = "  +1.35%  "

$ws.Range("D35").Value = "'3.661"
$ws.Range("D35").Style = $ws.Range("D6").Style
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "'0.2229"
$ws.Range("D36").Style = $ws.Range("D6").Style
$ws.Range("E36").Value = "  -2.29%  "

$ws.Range("D37").Value = "'0.02359"
$ws.Range("D37").Style = $ws.Range("D6").Style
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'5.197"
$ws.Range("D38").Style = $ws.Range("D6").Style
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("D39").Value = "'8.791"
$ws.Range("D39").Style = $ws.Range("D6").Style
$ws.Range("E39").Value = "  -1.67%  "

$ws.Range("D40").Value = "'0.6323"
$ws.Range("D40").Style = $ws.Range("D6").Style
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").Value = "'11.25"
$ws.Range("D41").Style = $ws.Range("D6").Style
$ws.Range("E41").Value = "  -1.80%  "

$ws.Range("D42").Value = "'1.195"
$ws.Range("D42").Style = $ws.Range("D6").Style
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("D43").Value = "'1.399"
$ws.Range("D43").Style = $ws.Range("D6").Style
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("D44").Value = "'13.55"
$ws.Range("D44").Style = $ws.Range("D6").Style
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "'3.778"
$ws.Range("D45").Style = $ws.Range("D6").Style
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("D46").Value = "'0.5941"
$ws.Range("D46").Style = $ws.Range("D6").Style
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").Value = "'126.03"
$ws.Range("D47").Style = $ws.Range("D6").Style
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").Value = "'1.987"
$ws.Range("D48").Style = $ws.Range("D6").Style
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").Value = "'1.201"
$ws.Range("D49").Style = $ws.Range("D6").Style
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").Value = "'74.32"
$ws.Range("D51").Style = $ws.Range("D6").Style
$ws.Range("E51").Value = "  -0.41%  "
